$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "303.76"
Set-TextValue $ws.Range("E2") "-0.06%"
Set-TextValue $ws.Range("D3") "35.69"
Set-TextValue $ws.Range("E3") "-3.95%"
Set-TextValue $ws.Range("D4") "5.060"
Set-TextValue $ws.Range("E4") "1.20%"
Set-TextValue $ws.Range("D5") "0.07862"
Set-TextValue $ws.Range("E5") "-0.03%"
Set-TextValue $ws.Range("D6") "2.095"
Set-TextValue $ws.Range("E6") "-5.48%"
Set-TextValue $ws.Range("D7") "7.919"
Set-TextValue $ws.Range("E7") "-1.07%"
Set-TextValue $ws.Range("D8") "0.9210"
Set-TextValue $ws.Range("E8") "0.14%"
Set-TextValue $ws.Range("D9") "0.09751"
Set-TextValue $ws.Range("E9") "1.73%"
Set-TextValue $ws.Range("D10") "0.1839"
Set-TextValue $ws.Range("E10") "-2.24%"
Set-TextValue $ws.Range("D11") "0.08663"
Set-TextValue $ws.Range("E11") "0.46%"
Set-TextValue $ws.Range("D12") "0.03558"
Set-TextValue $ws.Range("E12") "-1.44%"
Set-TextValue $ws.Range("D13") "0.09907"
Set-TextValue $ws.Range("E13") "-0.61%"
Set-TextValue $ws.Range("D14") "0.001438"
Set-TextValue $ws.Range("E14") "-3.07%"
Set-TextValue $ws.Range("D15") "0.005691"
Set-TextValue $ws.Range("E15") "-0.08%"
Set-TextValue $ws.Range("D16") "3.458"
Set-TextValue $ws.Range("E16") "0.16%"
Set-TextValue $ws.Range("D17") "4.127"
Set-TextValue $ws.Range("E17") "2.64%"
Set-TextValue $ws.Range("D18") "2.750"
Set-TextValue $ws.Range("E18") "21.98%"
Set-TextValue $ws.Range("D19") "0.3373"
Set-TextValue $ws.Range("E19") "-1.19%"
Set-TextValue $ws.Range("D20") "0.1341"
Set-TextValue $ws.Range("E20") "1.81%"
Set-TextValue $ws.Range("D21") "5.174"
Set-TextValue $ws.Range("E21") "8.65%"
Set-TextValue $ws.Range("D22") "0.2214"
Set-TextValue $ws.Range("E22") "0.61%"
Set-TextValue $ws.Range("D23") "0.04498"
Set-TextValue $ws.Range("E23") "-0.38%"
Set-TextValue $ws.Range("D24") "0.001237"
Set-TextValue $ws.Range("E24") "0.34%"
Set-TextValue $ws.Range("D25") "0.004859"
Set-TextValue $ws.Range("E25") "8.73%"
Set-TextValue $ws.Range("D26") "0.0001303"
Set-TextValue $ws.Range("E26") "-6.99%"
Set-TextValue $ws.Range("D27") "0.0004759"
Set-TextValue $ws.Range("E27") "0.18%"
Set-TextValue $ws.Range("D39") "0.01828"
Set-TextValue $ws.Range("E39") "-0.45%"
Set-TextValue $ws.Range("D40") "0.04697"
Set-TextValue $ws.Range("E40") "-1.02%"
Set-TextValue $ws.Range("D41") "0.007904"
Set-TextValue $ws.Range("E41") "-2.57%"
Set-TextValue $ws.Range("D42") "0.1388"
Set-TextValue $ws.Range("E42") "-0.54%"
Set-TextValue $ws.Range("D43") "0.007618"
Set-TextValue $ws.Range("E43") "0.94%"
Set-TextValue $ws.Range("D44") "0.002194"
Set-TextValue $ws.Range("E44") "-1.66%"
Set-TextValue $ws.Range("E45") "6.98%"
Set-TextValue $ws.Range("D46") "0.00006295"
Set-TextValue $ws.Range("E46") "-0.16%"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.22%"
Set-TextValue $ws.Range("E48") "0.18%"
Set-TextValue $ws.Range("D49") "50.57"
Set-TextValue $ws.Range("E49") "65.13%"
Set-TextValue $ws.Range("D50") "0.001904"
Set-TextValue $ws.Range("E50") "10.65%"
Set-TextValue $ws.Range("D51") "0.00002105"
Set-TextValue $ws.Range("E51") "0.22%"
